$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Replace the company names in column A (rows 2-7) with the new set.
$ws.Range("A2").Value = "AliBaba"
$ws.Range("A3").Value = "DetaCom"
$ws.Range("A4").Value = "Alphabet"
$ws.Range("A5").Value = "Gameboy"
$ws.Range("A6").Value = "Riot"
$ws.Range("A7").Value = "Steam"

# Remove the old rows 8 and 9 (Samsung / Lenovo) entirely.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(8).Delete()

# Update the selection to match the saved state.
$ws.Range("A8:B9").Select()

$wb.Save()
